$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6293
$ws.Range("C22").Value = 994
$ws.Range("D22").Value = 5831171
$ws.Range("E22").Value = 926.6122675989194
$ws.Range("F22").Value = 8.331898777758639
$ws.Range("G22").Value = 3.974895397489542
$ws.Range("H22").Value = 26.80979740603171
